# Internship_Oscar_Veldman.pptx — date bump (1/17/2017 -> 1/18/2017) on the
# slide master + every slide layout's "datetimeFigureOut" footer field, and
# the closing slide's "Questions or comments?" -> "Any questions or
# comments?" text tweak.

$p = $ppt.ActivePresentation

$oldDate = "1/17/2017"
$newDate = "1/18/2017"

# --- Slide master: update the Date placeholder if it still shows the old
#     cached field text.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout: same Date placeholder update.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Closing slide: tweak the question prompt copy.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "Questions or comments?") {
                $shape.TextFrame.TextRange.Text = "Any questions or comments?"
            }
        }
    }
}
